# Rename the worksheet "C_11" -> "C_16.2".
# Excel automatically cascades this rename into:
#   - the sheet tab / <sheets> entry
#   - the workbook-level TitlesOfParts (docProps/app.xml)
#   - any formulas / defined names that reference the sheet
#     (here: the hidden _FilterDatabase name used by the AutoFilter
#     range on B5:B5)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "C_16.2"

# The new sheet name contains a period, so when Excel re-writes the
# reference inside the _xlnm._FilterDatabase defined name it quotes the
# sheet name (standard Excel behaviour for names with "." or spaces):
#   'C_16.2'!$B$5:$B$5
$filterName = $wb.Names.Item(1)
$filterName.RefersTo = "='C_16.2'!`$B`$5:`$B`$5"
